# Apply the "Saldo" export update described by the commit diff.
# Strategy: walk the sheet top-to-bottom performing the same sequence of
# row delete / row insert / value update operations that the diff encodes,
# so each operation's row number already accounts for the shifts caused by
# the operations that precede it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 2 (005642649 / VR / 500129.35) was removed entirely.
$ws.Rows.Item(2).Delete()

# 2) Row 2 (now 004212438 / KENIA) balance updated.
$ws.Cells.Item(2, 3).Value = 290404.58

# 3) A new row (004459461 / INTERLAGOS / 100056.02) was inserted before the
#    004480970 / ALBERTO row (currently row 4).
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004459461"
$ws.Cells.Item(4, 2).Value = "INTERLAGOS"
$ws.Cells.Item(4, 3).Value = 100056.02

# 4) Row 5 (004480970 / ALBERTO) balance updated.
$ws.Cells.Item(5, 3).Value = 67645.49

# 5) A new row (004376145 / LUCYENE / 25836.21) was inserted before the
#    004321016 / JOAQUIM row (currently row 8).
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "004376145"
$ws.Cells.Item(8, 2).Value = "LUCYENE"
$ws.Cells.Item(8, 3).Value = 25836.21

# 6) Row 10 (002694089 / VITOR) was removed entirely.
$ws.Rows.Item(10).Delete()

# 7) Row 13 (004265173 / JULIA) balance updated.
$ws.Cells.Item(13, 3).Value = 5306.54

# 8) Row 14 (004805333 / ROSANA) balance updated.
$ws.Cells.Item(14, 3).Value = 3802.6

# 9) The old 004459461 / INTERLAGOS / 56.02 row (now row 79) was removed -
#    that account now appears earlier in the sheet with its new balance
#    (see step 3).
$ws.Rows.Item(79).Delete()

# 10) A new row (005366255 / RAPHAELA / 40.81) was inserted right after the
#     005000460 / MARIANA row (currently row 96), i.e. as row 97.
$ws.Rows.Item(97).Insert()
$ws.Cells.Item(97, 1).NumberFormat = "@"
$ws.Cells.Item(97, 1).Value = "005366255"
$ws.Cells.Item(97, 2).Value = "RAPHAELA"
$ws.Cells.Item(97, 3).Value = 40.81
